$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '79.626.60'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +4.34%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.194.06'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +5.29%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '205.20'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +2.58%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '637.54'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +2.27%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.246'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +20.14%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.606'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +10.34%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '3.189.97'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +5.20%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.634'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +44.23%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.166'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +3.26%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000241'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +25.21%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.45'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +3.81%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.769.35'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +4.87%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '32.16'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +10.65%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '79.165.72'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +3.77%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.193.06'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +5.50%  '
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +8.66%  '
$ws.Range('B20').NumberFormat = '@'
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').NumberFormat = '@'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '9.38'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +4.67%  '
$ws.Range('B21').NumberFormat = '@'
$ws.Range('B21').Value = 'SuiNetwork'
$ws.Range('C21').NumberFormat = '@'
$ws.Range('C21').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '3.00'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +30.73%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '435.69'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +16.47%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.22'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +20.38%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.92'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +12.90%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.355.14'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +4.96%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '77.45'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +5.71%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +13.47%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.997'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.06%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0000121'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +10.50%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.19'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +11.84%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.993'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.89%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.49'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +6.38%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '531.12'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +7.71%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.01'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +3.10%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.143'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +24.65%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '23.33'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +13.59%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.123'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +16.61%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.11%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.413'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +7.80%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '165.11'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +1.40%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '20.03'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.03%  '
$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '192.49'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +1.79%  '
$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'USDe'
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.00'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.02%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.56'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +9.29%  '
$ws.Range('B45').NumberFormat = '@'
$ws.Range('B45').Value = 'Stacks'
$ws.Range('C45').NumberFormat = '@'
$ws.Range('C45').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.81'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +10.01%  '
$ws.Range('B46').NumberFormat = '@'
$ws.Range('B46').Value = 'Mantle'
$ws.Range('C46').NumberFormat = '@'
$ws.Range('C46').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.806'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +1.19%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.34'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +6.06%  '
$ws.Range('B48').NumberFormat = '@'
$ws.Range('B48').Value = 'dogwifhat'
$ws.Range('C48').NumberFormat = '@'
$ws.Range('C48').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.63'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +7.24%  '
$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'OKB'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '43.56'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +3.73%  '
$ws.Range('B50').NumberFormat = '@'
$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').NumberFormat = '@'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '26.04'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +17.11%  '
$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value = 'ARBITRUM'
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.641'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +6.13%  '
